# Insert a new "weight" column between "bmi" (F) and "text_while_driving_30d"
# (previously G, now shifted to H). This pushes the existing G/H/I columns
# (text_while_driving_30d, smoked_ever, bullied_past_12mo) one column to the
# right, becoming H/I/J respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("G:G").Insert()

$ws.Range("G1").Value = "weight"

$weights = @{
    2  = 3.1945
    3  = 1.0206
    4  = 2.0312
    5  = 0.5954
    6  = 0.6833
    7  = 0.8456
    8  = 0.4695
    9  = 0.8511
    10 = 1.2463
    11 = 0.8134
    12 = 2.322
    13 = 0.2678
    14 = 0.8173
    15 = 0.4251
    16 = 0.5173
    17 = 1.1997
    18 = 0.9901
    19 = 4.4926
    20 = 0.8376
    21 = 0.2913
}

foreach ($row in $weights.Keys) {
    $ws.Cells.Item($row, 7).Value = $weights[$row]
}
